$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.133.39"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "1.652.69"
$ws.Range("E3").Value = "  -3.41%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'215.70"
$ws.Range("E5").Value = "  -3.89%  "
$ws.Range("D6").Value = "'0.5105"
$ws.Range("E6").Value = "  -3.48%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").Value = "'19.95"
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("D11").Value = "'0.07788"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "1.653.95"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "'4.279"
$ws.Range("E13").Value = "  -5.00%  "
$ws.Range("D14").Value = "1.879.40"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "'0.5517"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "'64.00"
$ws.Range("E17").Value = "  -5.96%  "
$ws.Range("D18").Value = "26.127.52"
$ws.Range("E18").Value = "  -4.51%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'210.02"
$ws.Range("E20").Value = "  -5.72%  "
$ws.Range("D21").Value = "'4.414"
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("E22").Value = "  -2.98%  "
$ws.Range("D23").Value = "'6.038"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'143.66"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").Value = "'1.740"
$ws.Range("E26").Value = "  +3.19%  "
$ws.Range("D27").Value = "'0.1179"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'6.971"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").Value = "'15.83"
$ws.Range("E29").Value = "  -2.65%  "
$ws.Range("D30").Value = "'0.05098"
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("D32").Value = "'3.342"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").Value = "'3.221"
$ws.Range("E33").Value = "  -6.27%  "
$ws.Range("D34").Value = "'1.567"
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("D35").Value = "'2.748"
$ws.Range("E35").Value = "  -4.44%  "
$ws.Range("D36").Value = "'2.363"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'0.9255"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").Value = "1.167.04"
$ws.Range("E38").Value = "  +5.81%  "
$ws.Range("D39").Value = "'0.5680"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").Value = "'0.01590"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").Value = "'2.549"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").Value = "'0.8321"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "'5.657"
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").Value = "'100.42"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "1.790.23"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").Value = "0.0₈116"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'0.4550"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").Value = "'55.65"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "'1.006"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'7.917"
$ws.Range("E51").Value = "  -2.18%  "
